# "Generate Report for handback"
#
# Fills in the handback columns (Latest Target File / Latest Handback File /
# Latest Handback DateTime / Handoff Reason) for the two already-handed-off
# rows on each language sheet, and flips the Status column from
# "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

function Get-HyperlinkAddress($ws, $addr) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            return $hl.Address
        }
    }
    return $null
}

function Update-LanguageSheet($sheetName, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Grab the existing hyperlink targets so the new "Latest Target File" /
    # "Latest Handback File" links point at the same files already used for
    # "Source File Name" / "Latest Handoff File".
    $mdUrl = Get-HyperlinkAddress $ws '$A$2'
    $xlfUrl = Get-HyperlinkAddress $ws '$C$2'
    $mdDisplay = $ws.Range("A2").Value2
    $xlfDisplay = $ws.Range("C2").Value2

    foreach ($row in 2, 3) {
        # Status: handed off -> handed back
        $ws.Range("B$row").Value = "Handed back: in sync with en-US"

        # Latest Target File (E) / Latest Handback File (F): same md/xlf
        # pair that was handed off, now coming back.
        $ws.Hyperlinks.Add($ws.Range("E$row"), $mdUrl, "", "", $mdDisplay) | Out-Null
        $ws.Range("E$row").Style = "HyperLink"

        $ws.Hyperlinks.Add($ws.Range("F$row"), $xlfUrl, "", "", $xlfDisplay) | Out-Null
        $ws.Range("F$row").Style = "HyperLink"

        # Latest Handback DateTime (G): was the "0001-01-01 00:00:00" sentinel.
        $ws.Range("G$row").Value = $handbackDateTime

        # Handoff Reason (H) stays "Include" for these two rows.
        $ws.Range("H$row").Value = "Include"
    }
}

Update-LanguageSheet "zh-cn" "2016-01-18 11:27:33"
Update-LanguageSheet "de-de" "2016-01-18 11:27:49"

# The "Overview" sheet mirrors the same Status string ("Ready for handoff")
# for both languages on these two rows; keep it in sync with the handback.
$overview = $wb.Worksheets.Item("Overview")
foreach ($row in 2, 3) {
    $overview.Range("B$row").Value = "Handed back: in sync with en-US"
    $overview.Range("C$row").Value = "Handed back: in sync with en-US"
}
